$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043260161776686
$ws.Range("D2").Value = 1.046428770128438
$ws.Range("E2").Value = 1.041230987379512
$ws.Range("F2").Value = 1.054716789940024
$ws.Range("I2").Value = 1.035695167092103
$ws.Range("J2").Value = 1.048331565347354
$ws.Range("K2").Value = 1.049194232275643
$ws.Range("L2").Value = 1.044011088834224
$ws.Range("M2").Value = 1.057459232750838
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044776125857301
$ws.Range("D3").Value = 1.047574752030882
$ws.Range("E3").Value = 1.042538705071644
$ws.Range("F3").Value = 1.056094196138116
$ws.Range("I3").Value = 1.035997814210862
$ws.Range("J3").Value = 1.049491494871726
$ws.Range("K3").Value = 1.050150986211601
$ws.Range("L3").Value = 1.045128086639669
$ws.Range("M3").Value = 1.058648496899758
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045755743955566
$ws.Range("D4").Value = 1.048314926298311
$ws.Range("E4").Value = 1.043383953425449
$ws.Range("F4").Value = 1.056984470582771
$ws.Range("I4").Value = 1.036191685584828
$ws.Range("J4").Value = 1.050240359924661
$ws.Range("K4").Value = 1.050768141607328
$ws.Range("L4").Value = 1.045849387169564
$ws.Range("M4").Value = 1.059416500872545
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046167268682975
$ws.Range("D5").Value = 1.048625776260939
$ws.Range("E5").Value = 1.043739078018836
$ws.Range("F5").Value = 1.057358508280055
$ws.Range("I5").Value = 1.036272721084111
$ws.Range("J5").Value = 1.050554784244984
$ws.Range("K5").Value = 1.051027136324884
$ws.Range("L5").Value = 1.046152274529214
$ws.Range("M5").Value = 1.0597390086726
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046236347723036
$ws.Range("D6").Value = 1.048677950744643
$ws.Range("E6").Value = 1.043798692384659
$ws.Range("F6").Value = 1.057421297314224
$ws.Range("I6").Value = 1.036286299901194
$ws.Range("J6").Value = 1.050607554254351
$ws.Range("K6").Value = 1.051070595977688
$ws.Range("L6").Value = 1.046203110422392
$ws.Range("M6").Value = 1.059793138125191
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04576124396519
$ws.Range("D7").Value = 1.04831908113977
$ws.Range("E7").Value = 1.04338869946682
$ws.Range("F7").Value = 1.056989469404074
$ws.Range("I7").Value = 1.036192770222227
$ws.Range("J7").Value = 1.050244562837297
$ws.Range("K7").Value = 1.050771604099502
$ws.Range("L7").Value = 1.045853435724151
$ws.Range("M7").Value = 1.059420811650039
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043772765054363
$ws.Range("D8").Value = 1.046816342888007
$ws.Range("E8").Value = 1.041673133030759
$ws.Range("F8").Value = 1.055182501138211
$ws.Range("I8").Value = 1.035797855204748
$ws.Range("J8").Value = 1.048723921447507
$ws.Range("K8").Value = 1.049517973595136
$ws.Range("L8").Value = 1.044388891197042
$ws.Range("M8").Value = 1.057861469891345
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040258430923222
$ws.Range("D9").Value = 1.044157762431225
$ws.Range("E9").Value = 1.038642700955888
$ws.Range("F9").Value = 1.051990485707187
$ws.Range("I9").Value = 1.035086870374437
$ws.Range("J9").Value = 1.046031198798408
$ws.Range("K9").Value = 1.047293956184707
$ws.Range("L9").Value = 1.041796681934694
$ws.Range("M9").Value = 1.05510176787364
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037908096796549
$ws.Range("D10").Value = 1.042377985299029
$ws.Range("E10").Value = 1.036617105066969
$ws.Range("F10").Value = 1.049856805095536
$ws.Range("I10").Value = 1.034602631456029
$ws.Range("J10").Value = 1.044226867822962
$ws.Range("K10").Value = 1.045800949793096
$ws.Range("L10").Value = 1.04006050024566
$ws.Range("M10").Value = 1.053253614676645
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036888506212694
$ws.Range("D11").Value = 1.041605504495641
$ws.Range("E11").Value = 1.035738663799029
$ws.Range("F11").Value = 1.048931473296448
$ws.Range("I11").Value = 1.034390496916137
$ws.Range("J11").Value = 1.043443319255452
$ws.Range("K11").Value = 1.045151951151662
$ws.Range("L11").Value = 1.039306736441617
$ws.Range("M11").Value = 1.052451290549655
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036509493311695
$ws.Range("D12").Value = 1.041318290781312
$ws.Range("E12").Value = 1.035412162834609
$ws.Range("F12").Value = 1.048587541471989
$ws.Range("I12").Value = 1.034311329702885
$ws.Range("J12").Value = 1.043151928106745
$ws.Range("K12").Value = 1.044910500682838
$ws.Range("L12").Value = 1.039026450522923
$ws.Range("M12").Value = 1.052152955579729
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036590806205341
$ws.Range("D13").Value = 1.041379911857991
$ws.Range("E13").Value = 1.035482207955182
$ws.Range("F13").Value = 1.048661326191789
$ws.Range("I13").Value = 1.03432832816017
$ws.Range("J13").Value = 1.043214448275358
$ws.Range("K13").Value = 1.044962310067075
$ws.Range("L13").Value = 1.039086586681659
$ws.Range("M13").Value = 1.05221696384837
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036857182888417
$ws.Range("D14").Value = 1.041581769064319
$ws.Range("E14").Value = 1.035711679425209
$ws.Range("F14").Value = 1.048903048361708
$ws.Range("I14").Value = 1.034383960509682
$ws.Range("J14").Value = 1.04341923987534
$ws.Range("K14").Value = 1.045132000634342
$ws.Range("L14").Value = 1.039283574167228
$ws.Range("M14").Value = 1.052426636572437
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037021267447815
$ws.Range("D15").Value = 1.041706102684561
$ws.Range("E15").Value = 1.035853036506353
$ws.Range("F15").Value = 1.04905195169026
$ws.Range("I15").Value = 1.034418188214209
$ws.Range("J15").Value = 1.043545372650118
$ws.Range("K15").Value = 1.045236501618964
$ws.Range("L15").Value = 1.03940490412691
$ws.Range("M15").Value = 1.052555780785692
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037975723276682
$ws.Range("D16").Value = 1.0424292133597
$ws.Range("E16").Value = 1.036675375414057
$ws.Range("F16").Value = 1.049918185514573
$ws.Range("I16").Value = 1.034616658211118
$ws.Range("J16").Value = 1.044278821156249
$ws.Range("K16").Value = 1.045843968211426
$ws.Range("L16").Value = 1.040110482706685
$ws.Range("M16").Value = 1.053306818311588
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038573918395484
$ws.Range("D17").Value = 1.042882308832459
$ws.Range("E17").Value = 1.037190842099732
$ws.Range("F17").Value = 1.050461162291121
$ws.Range("I17").Value = 1.034740494256976
$ws.Range("J17").Value = 1.044738283969022
$ws.Range("K17").Value = 1.046224338333319
$ws.Range("L17").Value = 1.040552537556155
$ws.Range("M17").Value = 1.053777367905091
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038922654671969
$ws.Range("D18").Value = 1.043146416002676
$ws.Range("E18").Value = 1.037491375791688
$ws.Range("F18").Value = 1.050777733569847
$ws.Range("I18").Value = 1.034812488907735
$ws.Range("J18").Value = 1.045006063114694
$ws.Range("K18").Value = 1.046445959296409
$ws.Range("L18").Value = 1.040810189500827
$ws.Range("M18").Value = 1.054051633038856
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039041534333226
$ws.Range("D19").Value = 1.043236440082558
$ws.Range("E19").Value = 1.037593828283781
$ws.Range("F19").Value = 1.050885653015753
$ws.Range("I19").Value = 1.034836997115402
$ws.Range("J19").Value = 1.045097332148983
$ws.Range("K19").Value = 1.046521485346587
$ws.Range("L19").Value = 1.040898009915098
$ws.Range("M19").Value = 1.054145116886217
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038509756488727
$ws.Range("D20").Value = 1.042833714196483
$ws.Range("E20").Value = 1.037135550825975
$ws.Range("F20").Value = 1.05040292033456
$ws.Range("I20").Value = 1.034727232331331
$ws.Range("J20").Value = 1.044689010527073
$ws.Range("K20").Value = 1.046183553339896
$ws.Range("L20").Value = 1.040505129072631
$ws.Range("M20").Value = 1.053726902962927
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036778749716053
$ws.Range("D21").Value = 1.041522334922031
$ws.Range("E21").Value = 1.035644111613808
$ws.Range("F21").Value = 1.048831873402676
$ws.Range("I21").Value = 1.034367588427152
$ws.Range("J21").Value = 1.043358943466767
$ws.Range("K21").Value = 1.045082041620507
$ws.Range("L21").Value = 1.039225574717338
$ws.Range("M21").Value = 1.052364901957234
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035688705641418
$ws.Range("D22").Value = 1.040696196176043
$ws.Range("E22").Value = 1.034705173374244
$ws.Range("F22").Value = 1.047842803657764
$ws.Range("I22").Value = 1.034139318521245
$ws.Range("J22").Value = 1.042520670499307
$ws.Range("K22").Value = 1.044387256137385
$ws.Range("L22").Value = 1.038419303394552
$ws.Range("M22").Value = 1.051506726034815
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036266721994798
$ws.Range("D23").Value = 1.041134303662949
$ws.Range("E23").Value = 1.035203039362897
$ws.Range("F23").Value = 1.048367252871477
$ws.Range("I23").Value = 1.034260532958473
$ws.Range("J23").Value = 1.042965247311842
$ws.Range("K23").Value = 1.044755787413525
$ws.Range("L23").Value = 1.038846892419302
$ws.Range("M23").Value = 1.051961837092159
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03853874905152
$ws.Range("D24").Value = 1.0428556725641
$ws.Range("E24").Value = 1.037160534972424
$ws.Range("F24").Value = 1.050429237794106
$ws.Range("I24").Value = 1.03473322555654
$ws.Range("J24").Value = 1.044711275747579
$ws.Range("K24").Value = 1.046201983072373
$ws.Range("L24").Value = 1.040526551516601
$ws.Range("M24").Value = 1.05374970651162
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041168248390039
$ws.Range("D25").Value = 1.044846350533813
$ws.Range("E25").Value = 1.039427050552589
$ws.Range("F25").Value = 1.052816671496796
$ws.Range("I25").Value = 1.03527247602607
$ws.Range("J25").Value = 1.046728925872744
$ws.Range("K25").Value = 1.047870718495767
$ws.Range("L25").Value = 1.042468224506277
$ws.Range("M25").Value = 1.055816665287792

Write-Output "Applied 380 kV case updates"
